$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
# try unfreeze and refreeze with a new top-left
$excel.ActiveWindow.FreezePanes = $false
$ws2.Range("D13").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("C58").Select()
